$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.986.16'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.017.94'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.99%  '
$ws.Range('E6').Value = '  -3.67%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.80'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.379'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0785'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  -3.88%  '
$ws.Range('D12').Value = '2.322.08'
$ws.Range('E12').Value = '  -2.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.30'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.741'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.66%  '
$ws.Range('D17').Value = '2.008.70'
$ws.Range('E17').Value = '  -2.60%  '
$ws.Range('D18').Value = '36.904.46'
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('D21').Value = '0.0₃0821'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('E24').Value = '  +2.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.96%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.54%  '
$ws.Range('E28').Value = '  -5.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.34'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.42%  '
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.15'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.28'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.12'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('D41').Value = '1.485.02'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0217'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.32%  '
$ws.Range('E44').Value = '  -3.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.13'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.36'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').Value = '2.211.62'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.54%  '
